$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "o aplicativo deve permitir o gerenciamento de ar-condicionados;"
#    -> "o aplicativo deve permitir o gerenciamento de aparelhos de ar-condicionado;"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*gerenciamento de ar-condicionados*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $r = $d.Range($start, $end)
        $r.Text = "o aplicativo deve permitir o gerenciamento de aparelhos de ar-condicionado;"
        break
    }
}

# ---------------------------------------------------------------------------
# 2) "o aplicativo deve permitir regulagem de temperatura do ar-condicionado;"
#    -> "o aplicativo deve permitir regulagem de temperatura do aparelho de ar-condicionado;"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*regulagem de temperatura do ar-condicionado*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $r = $d.Range($start, $end)
        $r.Text = "o aplicativo deve permitir regulagem de temperatura do aparelho de ar-condicionado;"
        break
    }
}

# ---------------------------------------------------------------------------
# 3) "o aplicativo deve permitir o acionamento e desligamento automático de ar condicionado com base em presença; "
#    -> "o aplicativo deve permitir o acionamento e desligamento automático de aparelho de ar-condicionado com base em presença; "
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*acionamento e desligamento autom*tico de ar*condicionado*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $r = $d.Range($start, $end)
        $r.Text = "o aplicativo deve permitir o acionamento e desligamento automático de aparelho de ar-condicionado com base em presença;$([char]0x00A0)"
        break
    }
}

# ---------------------------------------------------------------------------
# 4) "o aplicativo deve permitir o agendamento do acionamento/desligamento do ar-condicionado; "
#    -> "o aplicativo deve permitir o agendamento do acionamento/desligamento do aparelho do ar-condicionado; "
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*agendamento do acionamento/desligamento do ar-condicionado*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        $r = $d.Range($start, $end)
        $r.Text = "o aplicativo deve permitir o agendamento do acionamento/desligamento do aparelho do ar-condicionado;$([char]0x00A0)"
        break
    }
}
